# Update the "Förändrad" (Changed) date column (C) for all data rows
# from its previous value to the new date serial number 45177
# (2023-09-08), leaving all other cells/columns untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45177
}
